$wb = $excel.ActiveWorkbook

# --- Match Records sheet: add new playtest rows (7-16) ---
$ws = $wb.Worksheets.Item("Match Records")

# Copy the number formats from the existing data rows (2-6) down onto the
# previously-blank rows 7-16 so the new cells pick up the same styles
# (date format on A, month-year format on B) instead of minting new ones.
# Only columns A:B carry a non-default style in this sheet, so restrict the
# paste to those columns to avoid stamping stray empty cells into C:H.
$ws.Range("A2:B6").Copy()
$ws.Range("A7:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$rows = @(
    @{D="Chaos Warriors"; E=0; F=2; G=0; H="Lose"},
    @{D="Chaos Warriors"; E=2; F=0; G=0; H="Win"},
    @{D="Chaos Warriors"; E=0; F=2; G=0; H="Lose"},
    @{D="Chaos Warriors"; E=0; F=2; G=0; H="Lose"},
    @{D="Chaos Warriors"; E=1; F=2; G=0; H="Lose"},
    @{D="Warriors";       E=0; F=2; G=0; H="Lose"},
    @{D="Warriors";       E=2; F=0; G=0; H="Win"},
    @{D="Warriors";       E=1; F=2; G=1; H="Lose"},
    @{D="Warriors";       E=1; F=2; G=0; H="Lose"},
    @{D="Warriors";       E=0; F=0; G=0; H=$null}
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = 45626
    $ws.Cells.Item($r, 2).Value = 38565
    $ws.Cells.Item($r, 3).Value = "Goat"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
    $r++
}

$ws.Range("H16").Select()

# --- Decks sheet: selection moved, and it becomes the active tab ---
$decks = $wb.Worksheets.Item("Decks")
$decks.Activate()
$decks.Range("C11").Select()
